$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StudentData")

$newRows = @(
    @("Kilback Satterfield", "Adolfo", "5860732332", "5948266638", "KilbackSatterfield@yopmail.com", "AD9849"),
    @("Waters Beer", "Earline", "3635531881", "4454075079", "WatersBeer@yopmail.com", "AD7047"),
    @("Gottlieb Kessler", "Maritza", "3483563899", "3156645634", "GottliebKessler@yopmail.com", "AD8870"),
    @("Lang Bernhard", "Dale", "7664871913", "8189668962", "LangBernhard@yopmail.com", "AD3067"),
    @("Schumm Heidenreich", "Titus", "8963008534", "5046953430", "SchummHeidenreich@yopmail.com", "AD5644")
)

$startRow = 19
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$c]
    }
}
